$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.368.61"
$ws.Range("E2").Value = "  +1.73%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.099.16"
$ws.Range("E3").Value = "  +1.29%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.72"
$ws.Range("E5").Value = "  +1.99%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.00"
$ws.Range("E6").Value = "  +3.22%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.096.81"
$ws.Range("E8").Value = "  +1.35%  "

# Row 9
$ws.Range("E9").Value = "  +1.65%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  +2.79%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.20"
$ws.Range("E11").Value = "  -4.60%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.474"
$ws.Range("E12").Value = "  +4.01%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("E13").Value = "  +1.30%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.30"
$ws.Range("E14").Value = "  +1.56%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.602.53"
$ws.Range("E15").Value = "  +1.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.441.66"
$ws.Range("E16").Value = "  +1.82%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.099.68"
$ws.Range("E17").Value = "  +1.04%  "

# Row 18
$ws.Range("E18").Value = "  +1.44%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.78"
$ws.Range("E19").Value = "  +0.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "484.29"
$ws.Range("E20").Value = "  +0.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.09"
$ws.Range("E21").Value = "  +2.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.680"
$ws.Range("E22").Value = "  +0.69%  "

# Row 23
$ws.Range("E23").Value = "  +4.54%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.17"
$ws.Range("E24").Value = "  +11.55%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.60"
$ws.Range("E25").Value = "  +0.71%  "

# Row 26
$ws.Range("E26").Value = "  +0.11%  "

# Row 27
$ws.Range("E27").Value = "  +1.82%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.05"
$ws.Range("E28").Value = "  +2.53%  "

# Row 29
$ws.Range("E29").Value = "  +2.68%  "

# Row 31
$ws.Range("E31").Value = "  +1.24%  "

# Row 32
$ws.Range("E32").Value = "  -0.70%  "

# Row 33
$ws.Range("E33").Value = "  +1.63%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.66"
$ws.Range("E34").Value = "  +0.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.24"
$ws.Range("E35").Value = "  +4.36%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.67"
$ws.Range("E36").Value = "  +0.54%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0411"
$ws.Range("E37").Value = "  +3.45%  "

# Row 38
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.99"
$ws.Range("E38").Value = "  +16.93%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "442.70"
$ws.Range("E39").Value = "  -4.35%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0817"
$ws.Range("E40").Value = "  -0.61%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.971.15"
$ws.Range("E41").Value = "  -2.52%  "

# Row 42
$ws.Range("E42").Value = "  +0.38%  "

# Row 43
$ws.Range("E43").Value = "  -3.78%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.24"
$ws.Range("E44").Value = "  +1.65%  "

# Row 45
$ws.Range("E45").Value = "  +2.78%  "

# Row 46
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("E46").Value = "  +5.42%  "

# Row 47
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.02%  "

# Row 48
$ws.Range("E48").Value = "  +1.99%  "

# Row 49
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0520"
$ws.Range("E49").Value = "  +2.22%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "118.37"
$ws.Range("E50").Value = "  +1.31%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.11"
$ws.Range("E51").Value = "  +1.62%  "
